$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update computed NATMI LR-pair statistics for rows 2-17 following recount of
# ligand/receptor-expressing cells (Dr Hou advice).
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 31.29437
$ws.Cells.Item(2, 8).Value = 93.88310999999999
$ws.Cells.Item(2, 9).Value = 0.2388439922596655
$ws.Cells.Item(2, 10).Value = 0.2388439922596655
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 13.441269
$ws.Cells.Item(2, 14).Value = 40.323807
$ws.Cells.Item(2, 15).Value = 0.08973082133481231
$ws.Cells.Item(2, 16).Value = 0.08973082133481232
$ws.Cells.Item(2, 17).Value = 420.63604535553
$ws.Cells.Item(2, 18).Value = 3785.72440819977
$ws.Cells.Item(2, 19).Value = 0.02143166759634534
$ws.Cells.Item(2, 20).Value = 0.02143166759634535
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 31.29437
$ws.Cells.Item(3, 8).Value = 93.88310999999999
$ws.Cells.Item(3, 9).Value = 0.2388439922596655
$ws.Cells.Item(3, 10).Value = 0.2388439922596655
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 54.711535
$ws.Cells.Item(3, 14).Value = 164.134605
$ws.Cells.Item(3, 15).Value = 0.3652416280068742
$ws.Cells.Item(3, 16).Value = 0.3652416280068742
$ws.Cells.Item(3, 17).Value = 1712.16301955795
$ws.Cells.Item(3, 18).Value = 15409.46717602155
$ws.Cells.Item(3, 19).Value = 0.08723576857258149
$ws.Cells.Item(3, 20).Value = 0.0872357685725815
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 31.29437
$ws.Cells.Item(4, 8).Value = 93.88310999999999
$ws.Cells.Item(4, 9).Value = 0.2388439922596655
$ws.Cells.Item(4, 10).Value = 0.2388439922596655
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 63.67711
$ws.Cells.Item(4, 14).Value = 191.03133
$ws.Cells.Item(4, 15).Value = 0.4250937452800914
$ws.Cells.Item(4, 16).Value = 0.4250937452800915
$ws.Cells.Item(4, 17).Value = 1992.7350408707
$ws.Cells.Item(4, 18).Value = 17934.6153678363
$ws.Cells.Item(4, 19).Value = 0.1015310872073104
$ws.Cells.Item(4, 20).Value = 0.1015310872073104
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 31.29437
$ws.Cells.Item(5, 8).Value = 93.88310999999999
$ws.Cells.Item(5, 9).Value = 0.2388439922596655
$ws.Cells.Item(5, 10).Value = 0.2388439922596655
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 17.96553866666667
$ws.Cells.Item(5, 14).Value = 53.896616
$ws.Cells.Item(5, 15).Value = 0.119933805378222
$ws.Cells.Item(5, 16).Value = 0.119933805378222
$ws.Cells.Item(5, 17).Value = 562.2202142839733
$ws.Cells.Item(5, 18).Value = 5059.98192855576
$ws.Cells.Item(5, 19).Value = 0.02864546888342829
$ws.Cells.Item(5, 20).Value = 0.02864546888342829
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 33.86972933333333
$ws.Cells.Item(6, 8).Value = 101.609188
$ws.Cells.Item(6, 9).Value = 0.2584995758255442
$ws.Cells.Item(6, 10).Value = 0.2584995758255442
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 13.441269
$ws.Cells.Item(6, 14).Value = 40.323807
$ws.Cells.Item(6, 15).Value = 0.08973082133481231
$ws.Cells.Item(6, 16).Value = 0.08973082133481232
$ws.Cells.Item(6, 17).Value = 455.252142926524
$ws.Cells.Item(6, 18).Value = 4097.269286338716
$ws.Cells.Item(6, 19).Value = 0.02319537925352667
$ws.Cells.Item(6, 20).Value = 0.02319537925352667
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 33.86972933333333
$ws.Cells.Item(7, 8).Value = 101.609188
$ws.Cells.Item(7, 9).Value = 0.2584995758255442
$ws.Cells.Item(7, 10).Value = 0.2584995758255442
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 54.711535
$ws.Cells.Item(7, 14).Value = 164.134605
$ws.Cells.Item(7, 15).Value = 0.3652416280068742
$ws.Cells.Item(7, 16).Value = 0.3652416280068742
$ws.Cells.Item(7, 17).Value = 1853.064881861193
$ws.Cells.Item(7, 18).Value = 16677.58393675074
$ws.Cells.Item(7, 19).Value = 0.09441480591360817
$ws.Cells.Item(7, 20).Value = 0.09441480591360818
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 33.86972933333333
$ws.Cells.Item(8, 8).Value = 101.609188
$ws.Cells.Item(8, 9).Value = 0.2584995758255442
$ws.Cells.Item(8, 10).Value = 0.2584995758255442
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 63.67711
$ws.Cells.Item(8, 14).Value = 191.03133
$ws.Cells.Item(8, 15).Value = 0.4250937452800914
$ws.Cells.Item(8, 16).Value = 0.4250937452800915
$ws.Cells.Item(8, 17).Value = 2156.726480428893
$ws.Cells.Item(8, 18).Value = 19410.53832386004
$ws.Cells.Item(8, 19).Value = 0.1098865528409955
$ws.Cells.Item(8, 20).Value = 0.1098865528409956
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 33.86972933333333
$ws.Cells.Item(9, 8).Value = 101.609188
$ws.Cells.Item(9, 9).Value = 0.2584995758255442
$ws.Cells.Item(9, 10).Value = 0.2584995758255442
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 17.96553866666667
$ws.Cells.Item(9, 14).Value = 53.896616
$ws.Cells.Item(9, 15).Value = 0.119933805378222
$ws.Cells.Item(9, 16).Value = 0.119933805378222
$ws.Cells.Item(9, 17).Value = 608.4879319675342
$ws.Cells.Item(9, 18).Value = 5476.391387707808
$ws.Cells.Item(9, 19).Value = 0.03100283781741376
$ws.Cells.Item(9, 20).Value = 0.03100283781741376
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 10.670404
$ws.Cells.Item(10, 8).Value = 32.011212
$ws.Cells.Item(10, 9).Value = 0.08143835106389757
$ws.Cells.Item(10, 10).Value = 0.08143835106389757
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 13.441269
$ws.Cells.Item(10, 14).Value = 40.323807
$ws.Cells.Item(10, 15).Value = 0.08973082133481231
$ws.Cells.Item(10, 16).Value = 0.08973082133481232
$ws.Cells.Item(10, 17).Value = 143.423770502676
$ws.Cells.Item(10, 18).Value = 1290.813934524084
$ws.Cells.Item(10, 19).Value = 0.007307530129116314
$ws.Cells.Item(10, 20).Value = 0.007307530129116316
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 10.670404
$ws.Cells.Item(11, 8).Value = 32.011212
$ws.Cells.Item(11, 9).Value = 0.08143835106389757
$ws.Cells.Item(11, 10).Value = 0.08143835106389757
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 54.711535
$ws.Cells.Item(11, 14).Value = 164.134605
$ws.Cells.Item(11, 15).Value = 0.3652416280068742
$ws.Cells.Item(11, 16).Value = 0.3652416280068742
$ws.Cells.Item(11, 17).Value = 583.79418191014
$ws.Cells.Item(11, 18).Value = 5254.147637191259
$ws.Cells.Item(11, 19).Value = 0.0297446759247733
$ws.Cells.Item(11, 20).Value = 0.0297446759247733
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 10.670404
$ws.Cells.Item(12, 8).Value = 32.011212
$ws.Cells.Item(12, 9).Value = 0.08143835106389757
$ws.Cells.Item(12, 10).Value = 0.08143835106389757
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 63.67711
$ws.Cells.Item(12, 14).Value = 191.03133
$ws.Cells.Item(12, 15).Value = 0.4250937452800914
$ws.Cells.Item(12, 16).Value = 0.4250937452800915
$ws.Cells.Item(12, 17).Value = 679.46048925244
$ws.Cells.Item(12, 18).Value = 6115.14440327196
$ws.Cells.Item(12, 19).Value = 0.03461893366318713
$ws.Cells.Item(12, 20).Value = 0.03461893366318714
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 10.670404
$ws.Cells.Item(13, 8).Value = 32.011212
$ws.Cells.Item(13, 9).Value = 0.08143835106389757
$ws.Cells.Item(13, 10).Value = 0.08143835106389757
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 17.96553866666667
$ws.Cells.Item(13, 14).Value = 53.896616
$ws.Cells.Item(13, 15).Value = 0.119933805378222
$ws.Cells.Item(13, 16).Value = 0.119933805378222
$ws.Cells.Item(13, 17).Value = 191.6995556509547
$ws.Cells.Item(13, 18).Value = 1725.296000858592
$ws.Cells.Item(13, 19).Value = 0.00976721134682081
$ws.Cells.Item(13, 20).Value = 0.009767211346820811
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 55.18980966666666
$ws.Cells.Item(14, 8).Value = 165.569429
$ws.Cells.Item(14, 9).Value = 0.4212180808508926
$ws.Cells.Item(14, 10).Value = 0.4212180808508926
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 13.441269
$ws.Cells.Item(14, 14).Value = 40.323807
$ws.Cells.Item(14, 15).Value = 0.08973082133481231
$ws.Cells.Item(14, 16).Value = 0.08973082133481232
$ws.Cells.Item(14, 17).Value = 741.8210777884669
$ws.Cells.Item(14, 18).Value = 6676.389700096202
$ws.Cells.Item(14, 19).Value = 0.03779624435582397
$ws.Cells.Item(14, 20).Value = 0.03779624435582397
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 55.18980966666666
$ws.Cells.Item(15, 8).Value = 165.569429
$ws.Cells.Item(15, 9).Value = 0.4212180808508926
$ws.Cells.Item(15, 10).Value = 0.4212180808508926
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 54.711535
$ws.Cells.Item(15, 14).Value = 164.134605
$ws.Cells.Item(15, 15).Value = 0.3652416280068742
$ws.Cells.Item(15, 16).Value = 0.3652416280068742
$ws.Cells.Item(15, 17).Value = 3019.519203221171
$ws.Cells.Item(15, 18).Value = 27175.67282899054
$ws.Cells.Item(15, 19).Value = 0.1538463775959112
$ws.Cells.Item(15, 20).Value = 0.1538463775959112
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 55.18980966666666
$ws.Cells.Item(16, 8).Value = 165.569429
$ws.Cells.Item(16, 9).Value = 0.4212180808508926
$ws.Cells.Item(16, 10).Value = 0.4212180808508926
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 63.67711
$ws.Cells.Item(16, 14).Value = 191.03133
$ws.Cells.Item(16, 15).Value = 0.4250937452800914
$ws.Cells.Item(16, 16).Value = 0.4250937452800915
$ws.Cells.Item(16, 17).Value = 3514.327581023396
$ws.Cells.Item(16, 18).Value = 31628.94822921057
$ws.Cells.Item(16, 19).Value = 0.1790571715685983
$ws.Cells.Item(16, 20).Value = 0.1790571715685983
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 55.18980966666666
$ws.Cells.Item(17, 8).Value = 165.569429
$ws.Cells.Item(17, 9).Value = 0.4212180808508926
$ws.Cells.Item(17, 10).Value = 0.4212180808508926
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 17.96553866666667
$ws.Cells.Item(17, 14).Value = 53.896616
$ws.Cells.Item(17, 15).Value = 0.119933805378222
$ws.Cells.Item(17, 16).Value = 0.119933805378222
$ws.Cells.Item(17, 17).Value = 991.5146595724738
$ws.Cells.Item(17, 18).Value = 8923.631936152264
$ws.Cells.Item(17, 19).Value = 0.05051828733055914
$ws.Cells.Item(17, 20).Value = 0.05051828733055914
